$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3.560699333333333
$ws.Range("H2").Value = 10.682098
$ws.Range("I2").Value = 0.2516303646515017
$ws.Range("J2").Value = 0.2516303646515017
$ws.Range("O2").Value = 0.7426786721750401
$ws.Range("P2").Value = 0.7426786721750401
$ws.Range("Q2").Value = 2.285817048828445
$ws.Range("R2").Value = 20.572353439456
$ws.Range("S2").Value = 0.1868805050982985
$ws.Range("T2").Value = 0.1868805050982985

$ws.Range("G3").Value = 3.560699333333333
$ws.Range("H3").Value = 10.682098
$ws.Range("I3").Value = 0.2516303646515017
$ws.Range("J3").Value = 0.2516303646515017
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.2224236666666667
$ws.Range("N3").Value = 0.667271
$ws.Range("O3").Value = 0.2573213278249599
$ws.Range("P3").Value = 0.2573213278249599
$ws.Range("Q3").Value = 0.7919838016175555
$ws.Range("R3").Value = 7.127854214557999
$ws.Range("S3").Value = 0.06474985955320327
$ws.Range("T3").Value = 0.06474985955320327

$ws.Range("I4").Value = 0.2153092375010323
$ws.Range("J4").Value = 0.2153092375010323
$ws.Range("O4").Value = 0.7426786721750401
$ws.Range("P4").Value = 0.7426786721750401
$ws.Range("S4").Value = 0.159905578614287
$ws.Range("T4").Value = 0.1599055786142871

$ws.Range("I5").Value = 0.2153092375010323
$ws.Range("J5").Value = 0.2153092375010323
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.2224236666666667
$ws.Range("N5").Value = 0.667271
$ws.Range("O5").Value = 0.2573213278249599
$ws.Range("P5").Value = 0.2573213278249599
$ws.Range("Q5").Value = 0.6776663407677777
$ws.Range("R5").Value = 6.098997066909999
$ws.Range("S5").Value = 0.05540365888674528
$ws.Range("T5").Value = 0.05540365888674528

$ws.Range("G6").Value = 4.835201333333333
$ws.Range("H6").Value = 14.505604
$ws.Range("I6").Value = 0.3416978971743455
$ws.Range("J6").Value = 0.3416978971743456
$ws.Range("O6").Value = 0.7426786721750401
$ws.Range("P6").Value = 0.7426786721750401
$ws.Range("Q6").Value = 3.103992954076444
$ws.Range("R6").Value = 27.935936586688
$ws.Range("S6").Value = 0.2537717405584463
$ws.Range("T6").Value = 0.2537717405584464

$ws.Range("G7").Value = 4.835201333333333
$ws.Range("H7").Value = 14.505604
$ws.Range("I7").Value = 0.3416978971743455
$ws.Range("J7").Value = 0.3416978971743456
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.2224236666666667
$ws.Range("N7").Value = 0.667271
$ws.Range("O7").Value = 0.2573213278249599
$ws.Range("P7").Value = 0.2573213278249599
$ws.Range("Q7").Value = 1.075463209631556
$ws.Range("R7").Value = 9.679168886684
$ws.Range("S7").Value = 0.08792615661589918
$ws.Range("T7").Value = 0.0879261566158992

$ws.Range("G8").Value = 2.707878
$ws.Range("H8").Value = 8.123634000000001
$ws.Range("I8").Value = 0.1913625006731204
$ws.Range("J8").Value = 0.1913625006731204
$ws.Range("O8").Value = 0.7426786721750401
$ws.Range("P8").Value = 0.7426786721750401
$ws.Range("Q8").Value = 1.738342139872
$ws.Range("R8").Value = 15.645079258848
$ws.Range("S8").Value = 0.1421208479040083
$ws.Range("T8").Value = 0.1421208479040083

$ws.Range("G9").Value = 2.707878
$ws.Range("H9").Value = 8.123634000000001
$ws.Range("I9").Value = 0.1913625006731204
$ws.Range("J9").Value = 0.1913625006731204
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.2224236666666667
$ws.Range("N9").Value = 0.667271
$ws.Range("O9").Value = 0.2573213278249599
$ws.Range("P9").Value = 0.2573213278249599
$ws.Range("Q9").Value = 0.6022961536460001
$ws.Range("R9").Value = 5.420665382814001
$ws.Range("S9").Value = 0.04924165276911211
$ws.Range("T9").Value = 0.04924165276911211
